$d = $word.ActiveDocument

# The document contains 3 occurrences of the sentence fragment
# "...program through an exception...". Each needs "through" changed to
# "throw", while the run that holds the sentence is split into three runs
# (matching the target OOXML):
#
#   "<prefix> program "   |   "throw"   |   " an exception. <suffix>"
#
# A plain Find & Replace on the Range merges the edited text back into a
# single run (losing the 3-way split the target requires). Wrapping the
# replaced word in a transient Bookmark forces the run boundaries to be
# created at the bookmark's edges without adding any residual character
# formatting; deleting the bookmark afterwards leaves the runs split apart.

$searchPhrase = "program through an exception"
$targetWord = "through"
$replacementWord = "throw"
$bookmarkName = "tmpSplitMark"

$replacedCount = 0
$maxIterations = 20

for ($iteration = 0; $iteration -lt $maxIterations; $iteration++) {
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute($searchPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    # $searchRange now spans the matched phrase; locate "through" inside it.
    $offset = $searchRange.Text.IndexOf($targetWord)
    if ($offset -lt 0) {
        break
    }

    $wordStart = $searchRange.Start + $offset
    $wordEnd = $wordStart + $targetWord.Length
    $wordRange = $d.Range($wordStart, $wordEnd)

    if ($wordRange.Text -ne $targetWord) {
        break
    }

    $d.Bookmarks.Add($bookmarkName, $wordRange) | Out-Null
    $wordRange.Text = $replacementWord
    $d.Bookmarks($bookmarkName).Delete()

    $replacedCount = $replacedCount + 1
}

Write-Output "Replaced $replacedCount occurrence(s) of '$targetWord' -> '$replacementWord'."
